# edit.ps1
# Commit: "fixed railroad legend/fixed Area Office spelling error/added link to minute orders"
#
# Adds a new log row (row 7) to the "SPM Updates Performed" tracking sheet for the
# cbardash commit: railroad legend/popup fix, Area Office service update, and the
# new Minute Orders popup link.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 7

# Give the row its final height up front so the new content doesn't get auto-sized
# to something else, matching the tall (wrapped, multi-paragraph) row used for the
# long "Description of Updates" text.
$ws.Rows.Item($row).RowHeight = 150

$ws.Cells.Item($row, 1).Value = 42324                  # A7 - Date (11/16/2015)
$ws.Cells.Item($row, 2).Value = "cbardash"              # B7 - Editor
# C7 (Commit Comment) is intentionally left blank for this entry.
$ws.Cells.Item($row, 4).Value = "Fixed Railroad Legend and Popup to show text descriptions instead of codes.  Replaced railroad service with new service.  Replaced InfoTemplate with PopupTemplate, since InfoTemplates cannot interpret coded domain values.
Changed Area Office layer to point to new service.  Added District Name to the popup window.
Added a link in the popup for the Highway Designation layer that points to the Minute Orders page based on the minute order number in the popup.
"                                                        # D7 - Description of Updates
$ws.Cells.Item($row, 5).Value = "see SPM_TestScrip.docx" # E7 - QA/QC Testing Procedure
$ws.Cells.Item($row, 6).Value = "YES"                    # F7 - Editor QC

# Re-assert the row height in case setting the long, wrapped text re-triggered an
# automatic row resize.
$ws.Rows.Item($row).RowHeight = 150

# Leave the selection where the author's edit session ended up.
$ws.Range("C7").Select()
